$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing "Quiz Dashboard Classes" result column (F17:F22) to include
# the newly added Class 11 / Class 12 sections.
$quizClasses = "Class 3-A, Class 3-B, Class 3-C, Class 4-A, Class 4-B, Class 4-C, Class 5-A, Class 5-B, Class 5-C, Class 6-A, Class 6-B, Class 6-C, Class 7-A, Class 7-B, Class 7-C, Class 8-A, Class 8-B, Class 8-C, Class 9-A, Class 9-B, Class 9-C, Class 10-A, Class 10-B, Class 10-C, Class 11-A, Class 11-B, Class 11-C, Class 12-A, Class 12-B, Class 12-C"
$ws.Range("F17").Value = $quizClasses
$ws.Range("F18").Value = $quizClasses
$ws.Range("F19").Value = $quizClasses
$ws.Range("F20").Value = $quizClasses
$ws.Range("F21").Value = $quizClasses
$ws.Range("F22").Value = $quizClasses

# Add new "searchAndViewContentSchool" test cases (rows 26-34)
$ws.Range("B26").Value = "Web"
$ws.Range("C26").Value = "searchAndViewContentSchool"
$ws.Range("D26").Value = "Admin"
$ws.Range("E26").Value = "Search String"
$ws.Range("F26").Value = "Autotrophic Nutrition"

$ws.Range("B27").Value = "Android"
$ws.Range("C27").Value = "searchAndViewContentSchool"
$ws.Range("D27").Value = "Admin"
$ws.Range("E27").Value = "Search String"
$ws.Range("F27").Value = "Autotrophic Nutrition"

$ws.Range("B28").Value = "iOS"
$ws.Range("C28").Value = "searchAndViewContentSchool"
$ws.Range("D28").Value = "Admin"
$ws.Range("E28").Value = "Search String"
$ws.Range("F28").Value = "Autotrophic Nutrition"

$ws.Range("B29").Value = "Web"
$ws.Range("C29").Value = "searchAndViewContentSchool"
$ws.Range("D29").Value = "Principal"
$ws.Range("E29").Value = "Search String"
$ws.Range("F29").Value = "Autotrophic Nutrition"

$ws.Range("B30").Value = "Android"
$ws.Range("C30").Value = "searchAndViewContentSchool"
$ws.Range("D30").Value = "Principal"
$ws.Range("E30").Value = "Search String"
$ws.Range("F30").Value = "Autotrophic Nutrition"

$ws.Range("B31").Value = "iOS"
$ws.Range("C31").Value = "searchAndViewContentSchool"
$ws.Range("D31").Value = "Principal"
$ws.Range("E31").Value = "Search String"
$ws.Range("F31").Value = "Autotrophic Nutrition"

$ws.Range("B32").Value = "Web"
$ws.Range("C32").Value = "searchAndViewContentSchool"
$ws.Range("D32").Value = "Teacher"
$ws.Range("E32").Value = "Search String"
$ws.Range("F32").Value = "Autotrophic Nutrition"

$ws.Range("B33").Value = "Android"
$ws.Range("C33").Value = "searchAndViewContentSchool"
$ws.Range("D33").Value = "Teacher"
$ws.Range("E33").Value = "Search String"
$ws.Range("F33").Value = "Autotrophic Nutrition"

$ws.Range("B34").Value = "iOS"
$ws.Range("C34").Value = "searchAndViewContentSchool"
$ws.Range("D34").Value = "Teacher"
$ws.Range("E34").Value = "Search String"
$ws.Range("F34").Value = "Autotrophic Nutrition"

# Add new "searchAndViewContentStudent" test cases (rows 35-43)
$ws.Range("B35").Value = "Web"
$ws.Range("C35").Value = "searchAndViewContentStudent"
$ws.Range("D35").Value = "Parent"
$ws.Range("E35").Value = "Search String"
$ws.Range("F35").Value = "Making stone tools and the discovery of fire"

$ws.Range("B36").Value = "Android"
$ws.Range("C36").Value = "searchAndViewContentStudent"
$ws.Range("D36").Value = "Parent"
$ws.Range("E36").Value = "Search String"
$ws.Range("F36").Value = "Making stone tools and the discovery of fire"

$ws.Range("B37").Value = "iOS"
$ws.Range("C37").Value = "searchAndViewContentStudent"
$ws.Range("D37").Value = "Parent"
$ws.Range("E37").Value = "Search String"
$ws.Range("F37").Value = "Making stone tools and the discovery of fire"

$ws.Range("B38").Value = "Web"
$ws.Range("C38").Value = "searchAndViewContentStudent"
$ws.Range("D38").Value = "Student"
$ws.Range("E38").Value = "Search String"
$ws.Range("F38").Value = "Making stone tools and the discovery of fire"

$ws.Range("B39").Value = "Android"
$ws.Range("C39").Value = "searchAndViewContentStudent"
$ws.Range("D39").Value = "Student"
$ws.Range("E39").Value = "Search String"
$ws.Range("F39").Value = "Making stone tools and the discovery of fire"

$ws.Range("B40").Value = "iOS"
$ws.Range("C40").Value = "searchAndViewContentStudent"
$ws.Range("D40").Value = "Student"
$ws.Range("E40").Value = "Search String"
$ws.Range("F40").Value = "Making stone tools and the discovery of fire"

$ws.Range("B41").Value = "Web"
$ws.Range("C41").Value = "searchAndViewContentStudent"
$ws.Range("D41").Value = "Guest"
$ws.Range("E41").Value = "Search String"
$ws.Range("F41").Value = "Making stone tools and the discovery of fire"

$ws.Range("B42").Value = "Android"
$ws.Range("C42").Value = "searchAndViewContentStudent"
$ws.Range("D42").Value = "Guest"
$ws.Range("E42").Value = "Search String"
$ws.Range("F42").Value = "Making stone tools and the discovery of fire"

$ws.Range("B43").Value = "iOS"
$ws.Range("C43").Value = "searchAndViewContentStudent"
$ws.Range("D43").Value = "Guest"
$ws.Range("E43").Value = "Search String"
$ws.Range("F43").Value = "Making stone tools and the discovery of fire"

# Match final selection / active cell from the authored edit
$ws.Range("F34").Select()
